$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2 through 43 from
# serial date 45753 (2025-04-06) to 45754 (2025-04-07).
$ws.Range("C2:C43").Value = 45754
